$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph ("Play Diamond Chance Slot for Free | Retro-Style Slot
#    Game"). The new paragraph has a leading empty run, a bold
#    "Meta description" run, and a plain run with the rest of the text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Diamond Chance Slot for Free | Retro-Style Slot Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Diamond Chance Slot and play for free. Enjoy simple gameplay and good chances of winning in this retro-style slot game.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titleRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicate bold
#    "Play Diamond Chance Slot for Free | Retro-Style Slot Game"
#    paragraph and replace the text of the following italic paragraph
#    with the new image-prompt text (keeping its leading empty run and
#    italic formatting).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)
$italicPara = $d.Paragraphs.Item($count)

$tailStart = $boldPara.Range.Start
$tailEnd = $italicPara.Range.End
$tailRange = $d.Range($tailStart, $tailEnd)

$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Please create a feature image for the game &quot;Diamond Chance&quot;. The image should be in cartoon style and feature a happy Maya warrior with glasses.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($tailXml)
